$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("testcases")
$ws2 = $wb.Worksheets.Item("data")

# --- Value changes ---
# testcases sheet: enable SearchTest (Runmode Y)
$ws1.Range("B4").Value = "Y"

# data sheet: RegisterTest table - toggle Runmode to Y for Firefox rows
$ws2.Range("A9").Value = "Y"
$ws2.Range("A11").Value = "Y"
# RegisterTest row 12 browser Chrome -> Edge
$ws2.Range("B12").Value = "Edge"

# SearchTest table: set runmode Y for all rows, set expected products, results
$ws2.Range("A16").Value = "Y"
$ws2.Range("A17").Value = "Y"
$ws2.Range("A18").Value = "Y"
$ws2.Range("D16").Value = "HP LP3065"
$ws2.Range("D17").Value = "There is no product that matches the search criteria."
$ws2.Range("D18").Value = "There is no product that matches the search criteria."
$ws2.Range("E17").Value = "Failure"
$ws2.Range("E18").Value = "Failure"
